$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to make room for the newest
# two reporting quarters (existing quarterly columns shift right by 2).
$ws.Columns("D:E").Insert()

# Copy number formats/styles from column F (which now holds what used
# to be column D) into the two new columns so the new cells inherit the
# correct date/number formatting used throughout the sheet.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D = most recent quarter,
# E = the quarter before it) with the latest financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 5900
$ws.Range("E8").Value = 5700
$ws.Range("D9:E9").Value = "NA"
$ws.Range("D10:E10").Value = "NA"
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14:E14").Value = 0
$ws.Range("D15:E15").Value = -200
$ws.Range("D17").Value = 600
$ws.Range("E17").Value = 500
$ws.Range("D18").Value = 5300
$ws.Range("E18").Value = 5200
$ws.Range("D20").Value = -3800
$ws.Range("E20").Value = -4500
$ws.Range("D21").Value = 1700
$ws.Range("E21").Value = 1000
$ws.Range("D22:E22").Value = 0
$ws.Range("D23").Value = 1500
$ws.Range("E23").Value = 800
$ws.Range("D24").Value = 300
$ws.Range("E24").Value = 100
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 1100
$ws.Range("E26").Value = 600
$ws.Range("D27").Value = 1000
$ws.Range("E27").Value = 500
$ws.Range("D28:E28").Value = 0
$ws.Range("D29:E29").Value = "NA"
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32").Value = 3800
$ws.Range("E32").Value = 4500
$ws.Range("D33").Value = 1000
$ws.Range("E33").Value = 500
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 1000
$ws.Range("E35").Value = 500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 4500
$ws.Range("E41").Value = 6600
$ws.Range("D42").Value = 114600
$ws.Range("E42").Value = 111200
$ws.Range("D43:E43").Value = 0
$ws.Range("D44:E44").Value = 0
$ws.Range("D45:E45").Value = 0
$ws.Range("D46:E46").Value = 0
$ws.Range("D47:E47").Value = 0
$ws.Range("D48:E48").Value = 14800
$ws.Range("D49:E49").Value = 0
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52:E52").Value = 0
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 632300
$ws.Range("E54").Value = 637400
$ws.Range("D57:E57").Value = 0
$ws.Range("D58:E58").Value = 0
$ws.Range("D59:E59").Value = 0
$ws.Range("D60:E60").Value = 0
$ws.Range("D61").Value = 10000
$ws.Range("E61").Value = 9900
$ws.Range("D62:E62").Value = 0
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 597800
$ws.Range("E66").Value = 604300
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = 14400
$ws.Range("E72").Value = 14200
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 34500
$ws.Range("E76").Value = 33100
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 1000
$ws.Range("E81").Value = 500
$ws.Range("D83:E83").Value = 300
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = -1700
$ws.Range("E89").Value = 4800
$ws.Range("D91").Value = 100
$ws.Range("E91").Value = -100
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = 9900
$ws.Range("E94").Value = -9500
$ws.Range("D96:E96").Value = 0
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -7000
$ws.Range("E100").Value = 22500
$ws.Range("D101:E101").Value = 0
$ws.Range("D102").Value = 1200
$ws.Range("E102").Value = 17700

Write-Host "Applied quarterly financial update"
